$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.640.15"
$ws.Range("E2").Value = "  +0.75%  "

$ws.Range("D3").Value = "3.398.91"
$ws.Range("E3").Value = "  -0.28%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("E5").Value = "  -0.92%  "

$ws.Range("D6").Value = "176.41"
$ws.Range("E6").Value = "  -0.97%  "

$ws.Range("D7").Value = "0.633"
$ws.Range("E7").Value = "  +0.51%  "

$ws.Range("D8").Value = "3.392.02"
$ws.Range("E8").Value = "  -0.11%  "

$ws.Range("E9").Value = "  +0.02%  "

$ws.Range("D10").Value = "0.172"
$ws.Range("E10").Value = "  +2.97%  "

$ws.Range("D11").Value = "0.641"
$ws.Range("E11").Value = "  +0.50%  "

$ws.Range("D12").Value = "53.74"
$ws.Range("E12").Value = "  -2.55%  "

$ws.Range("D13").Value = "0.0000278"
$ws.Range("E13").Value = "  -0.59%  "

$ws.Range("D14").Value = "9.22"
$ws.Range("E14").Value = "  +0.36%  "

$ws.Range("D15").Value = "3.941.15"
$ws.Range("E15").Value = "  -0.02%  "

$ws.Range("D16").Value = "18.34"
$ws.Range("E16").Value = "  -0.29%  "

$ws.Range("D17").Value = "3.422.62"
$ws.Range("E17").Value = "  +0.39%  "

$ws.Range("E18").Value = "  +0.57%  "

$ws.Range("D19").Value = "65.518.12"
$ws.Range("E19").Value = "  +0.59%  "

$ws.Range("D20").Value = "11.88"
$ws.Range("E20").Value = "  -0.73%  "

$ws.Range("E21").Value = "  +0.69%  "

$ws.Range("D22").Value = "482.06"
$ws.Range("E22").Value = "  +2.62%  "

$ws.Range("D23").Value = "4.94"
$ws.Range("E23").Value = "  -1.44%  "

$ws.Range("E24").Value = "  -1.05%  "

$ws.Range("D25").Value = "14.32"
$ws.Range("E25").Value = "  +5.02%  "

$ws.Range("D26").Value = "89.10"
$ws.Range("E26").Value = "  +2.99%  "

$ws.Range("D27").Value = "2.92"
$ws.Range("E27").Value = "  +1.26%  "

$ws.Range("E28").Value = "  -2.37%  "

$ws.Range("B29").Value = "Filecoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D29").Value = "8.75"
$ws.Range("E29").Value = "  -2.14%  "

$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").Value = "31.48"
$ws.Range("E30").Value = "  +2.78%  "

$ws.Range("D31").Value = "6.56"
$ws.Range("E31").Value = "  -2.70%  "

$ws.Range("D32").Value = "11.54"
$ws.Range("E32").Value = "  -0.56%  "

$ws.Range("D33").Value = "62.52"
$ws.Range("E33").Value = "  +3.69%  "

$ws.Range("D34").Value = "576.48"
$ws.Range("E34").Value = "  -1.73%  "

$ws.Range("E36").Value = "  -0.04%  "

$ws.Range("E37").Value = "  +3.95%  "

$ws.Range("E38").Value = "  -0.46%  "

$ws.Range("D39").Value = "36.05"
$ws.Range("E39").Value = "  -0.45%  "

$ws.Range("E40").Value = "  -0.27%  "

$ws.Range("D41").Value = "0.0₃0740"
$ws.Range("E41").Value = "  -3.63%  "

$ws.Range("D42").Value = "3.117.94"
$ws.Range("E42").Value = "  +0.00%  "

$ws.Range("B43").Value = "ThetaToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D43").Value = "2.80"
$ws.Range("E43").Value = "  -2.99%  "

$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "0.0418"
$ws.Range("E44").Value = "  +0.66%  "

$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").Value = "2.45"
$ws.Range("E46").Value = "  -3.63%  "

$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").Value = "3.16"
$ws.Range("E47").Value = "  -1.89%  "

$ws.Range("D48").Value = "0.999"
$ws.Range("E48").Value = "  -0.09%  "

$ws.Range("D49").Value = "140.36"
$ws.Range("E49").Value = "  +2.25%  "

$ws.Range("E50").Value = "  +0.05%  "

$ws.Range("D51").Value = "8.44"
$ws.Range("E51").Value = "  -0.62%  "
